$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column C ("Pre Experimental Phase") values for rows 2-23
$cValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 0
    11 = 3
    12 = 0
    13 = 0
    14 = 4
    15 = 5
    16 = 4
    17 = 4
    18 = 0
    19 = 2
    20 = 0
    21 = 5
    22 = 5
    23 = 1
}

foreach ($row in $cValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $cValues[$row]
}

# Total row: sum formula for column C
$ws.Range("C24").Formula = "=SUM(C2:C23)"

# Update the active selection to match the saved view state
$ws.Range("C22").Select()
